# Auto-generated edit script replicating the Typhon_Profits.xlsx commit
# "chore: update Sheets via scheduled runner".
# Refreshes the leve-profit calculation columns (H:N) for 32 rows across
# all 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect
# updated market-board prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 13).ClearContents()
$ws.Cells.Item(40, 8).Value = 2359.4
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 2359.4
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 2359.4
$ws.Cells.Item(40, 14).Value = -2709.4

$ws.Cells.Item(74, 8).Value = 6253855.5
$ws.Cells.Item(74, 9).Value = 3309.25
$ws.Cells.Item(74, 11).Value = 3309.25
$ws.Cells.Item(74, 13).Value = -2373.25

$ws.Cells.Item(77, 8).Value = 6253855.5
$ws.Cells.Item(77, 9).Value = 3309.25
$ws.Cells.Item(77, 11).Value = 16546.25
$ws.Cells.Item(77, 13).Value = -11866.25

$ws.Cells.Item(80, 8).Value = 12825992
$ws.Cells.Item(80, 10).Value = 16245723
$ws.Cells.Item(80, 12).Value = 48737169
$ws.Cells.Item(80, 14).Value = -48739165

$ws.Cells.Item(83, 8).Value = 12825992
$ws.Cells.Item(83, 10).Value = 16245723
$ws.Cells.Item(83, 12).Value = 146211507
$ws.Cells.Item(83, 14).Value = -146221491

$ws.Cells.Item(100, 8).Value = 2414
$ws.Cells.Item(100, 9).Value = 1999.6666
$ws.Cells.Item(100, 10).Value = 2724.75
$ws.Cells.Item(100, 11).Value = 1999.6666
$ws.Cells.Item(100, 12).Value = 2724.75
$ws.Cells.Item(100, 13).Value = -1458.6666
$ws.Cells.Item(100, 14).Value = -3806.75

$ws.Cells.Item(101, 8).Value = 441
$ws.Cells.Item(101, 9).Value = 255
$ws.Cells.Item(101, 10).Value = 1185
$ws.Cells.Item(101, 11).Value = 765
$ws.Cells.Item(101, 12).Value = 3555
$ws.Cells.Item(101, 13).Value = 857
$ws.Cells.Item(101, 14).Value = -6799

$ws.Cells.Item(129, 8).Value = 222963.02
$ws.Cells.Item(129, 10).Value = 244681.61
$ws.Cells.Item(129, 12).Value = 734044.83
$ws.Cells.Item(129, 14).Value = -744044.83

$ws.Cells.Item(132, 8).Value = 3176.4517
$ws.Cells.Item(132, 9).Value = 3518.077
$ws.Cells.Item(132, 10).Value = 1400
$ws.Cells.Item(132, 11).Value = 10554.231
$ws.Cells.Item(132, 12).Value = 4200
$ws.Cells.Item(132, 13).Value = -8024.231
$ws.Cells.Item(132, 14).Value = -9260

$ws.Cells.Item(137, 8).Value = 1125.7812
$ws.Cells.Item(137, 9).Value = 1096.3334
$ws.Cells.Item(137, 10).Value = 1182
$ws.Cells.Item(137, 11).Value = 3289.0002
$ws.Cells.Item(137, 12).Value = 3546
$ws.Cells.Item(137, 13).Value = -739.0001999999999
$ws.Cells.Item(137, 14).Value = -8646

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3511.08
$ws.Cells.Item(45, 9).Value = 3171.3572
$ws.Cells.Item(45, 10).Value = 3943.4546
$ws.Cells.Item(45, 11).Value = 3171.3572
$ws.Cells.Item(45, 12).Value = 3943.4546
$ws.Cells.Item(45, 13).Value = -2794.3572
$ws.Cells.Item(45, 14).Value = -4697.4546

$ws.Cells.Item(61, 8).Value = 3685.7083
$ws.Cells.Item(61, 9).Value = 3814.2778
$ws.Cells.Item(61, 11).Value = 3814.2778
$ws.Cells.Item(61, 13).Value = -3602.2778

$ws.Cells.Item(110, 8).Value = 919.26666
$ws.Cells.Item(110, 9).Value = 825.6923
$ws.Cells.Item(110, 10).Value = 1527.5
$ws.Cells.Item(110, 11).Value = 825.6923
$ws.Cells.Item(110, 12).Value = 1527.5
$ws.Cells.Item(110, 13).Value = 1219.3077
$ws.Cells.Item(110, 14).Value = -5617.5

$ws.Cells.Item(132, 8).Value = 12902.155
$ws.Cells.Item(132, 9).Value = 1451.4
$ws.Cells.Item(132, 10).Value = 52979.8
$ws.Cells.Item(132, 11).Value = 4354.200000000001
$ws.Cells.Item(132, 12).Value = 158939.4
$ws.Cells.Item(132, 13).Value = -1824.200000000001
$ws.Cells.Item(132, 14).Value = -163999.4

$ws.Cells.Item(136, 8).Value = 3685.7083
$ws.Cells.Item(136, 9).Value = 3814.2778
$ws.Cells.Item(136, 11).Value = 11442.8334
$ws.Cells.Item(136, 13).Value = -8892.8334

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 24374.834
$ws.Cells.Item(82, 10).Value = 40411.668
$ws.Cells.Item(82, 12).Value = 40411.668
$ws.Cells.Item(82, 14).Value = -41177.668

$ws.Cells.Item(85, 8).Value = 24374.834
$ws.Cells.Item(85, 10).Value = 40411.668
$ws.Cells.Item(85, 12).Value = 40411.668
$ws.Cells.Item(85, 14).Value = -43063.668

$ws.Cells.Item(86, 8).Value = 1733.4
$ws.Cells.Item(86, 9).Value = 1589.6
$ws.Cells.Item(86, 10).Value = 2092.9
$ws.Cells.Item(86, 11).Value = 1589.6
$ws.Cells.Item(86, 12).Value = 2092.9
$ws.Cells.Item(86, 13).Value = -466.5999999999999
$ws.Cells.Item(86, 14).Value = -4338.9

$ws.Cells.Item(89, 8).Value = 1733.4
$ws.Cells.Item(89, 9).Value = 1589.6
$ws.Cells.Item(89, 10).Value = 2092.9
$ws.Cells.Item(89, 11).Value = 7948
$ws.Cells.Item(89, 12).Value = 10464.5
$ws.Cells.Item(89, 13).Value = -2332
$ws.Cells.Item(89, 14).Value = -21696.5

$ws.Cells.Item(134, 8).Value = 5458.1724
$ws.Cells.Item(134, 9).Value = 5576.0454
$ws.Cells.Item(134, 10).Value = 5087.7144
$ws.Cells.Item(134, 11).Value = 16728.1362
$ws.Cells.Item(134, 12).Value = 15263.1432
$ws.Cells.Item(134, 13).Value = -14193.1362
$ws.Cells.Item(134, 14).Value = -20333.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(59, 8).Value = 22450
$ws.Cells.Item(59, 10).Value = 22450
$ws.Cells.Item(59, 12).Value = 22450
$ws.Cells.Item(59, 14).Value = -24740

$ws.Cells.Item(122, 8).Value = 1427.76
$ws.Cells.Item(122, 9).Value = 1155.875
$ws.Cells.Item(122, 10).Value = 1911.1111
$ws.Cells.Item(122, 11).Value = 3467.625
$ws.Cells.Item(122, 12).Value = 5733.3333
$ws.Cells.Item(122, 13).Value = -1017.625
$ws.Cells.Item(122, 14).Value = -10633.3333

$ws.Cells.Item(132, 8).Value = 2428.7144
$ws.Cells.Item(132, 9).Value = 1816.2858
$ws.Cells.Item(132, 11).Value = 5448.857400000001
$ws.Cells.Item(132, 13).Value = -2918.857400000001

$ws.Cells.Item(134, 8).Value = 847.3333
$ws.Cells.Item(134, 9).Value = 655.2222
$ws.Cells.Item(134, 11).Value = 1965.6666
$ws.Cells.Item(134, 13).Value = 569.3334

$ws.Cells.Item(141, 8).Value = 13115.091
$ws.Cells.Item(141, 10).Value = 13115.091
$ws.Cells.Item(141, 12).Value = 13115.091
$ws.Cells.Item(141, 14).Value = -23475.091

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 720.21
$ws.Cells.Item(131, 10).Value = 742.914
$ws.Cells.Item(131, 12).Value = 2228.742
$ws.Cells.Item(131, 14).Value = -12308.742

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1621.6897
$ws.Cells.Item(102, 9).Value = 1326.64
$ws.Cells.Item(102, 10).Value = 3465.75
$ws.Cells.Item(102, 11).Value = 1326.64
$ws.Cells.Item(102, 12).Value = 3465.75
$ws.Cells.Item(102, 13).Value = 295.3599999999999
$ws.Cells.Item(102, 14).Value = -6709.75

$ws.Cells.Item(113, 8).Value = 5708.56
$ws.Cells.Item(113, 9).Value = 6762.5
$ws.Cells.Item(113, 11).Value = 6762.5
$ws.Cells.Item(113, 13).Value = -4592.5

$ws.Cells.Item(132, 8).Value = 23157.115
$ws.Cells.Item(132, 9).Value = 4477.4443
$ws.Cells.Item(132, 11).Value = 13432.3329
$ws.Cells.Item(132, 13).Value = -10902.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 856193.7
$ws.Cells.Item(122, 9).Value = 1092397.5
$ws.Cells.Item(122, 11).Value = 3277192.5
$ws.Cells.Item(122, 13).Value = -3274742.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(29, 8).Value = 11999
$ws.Cells.Item(29, 9).Value = 11999
$ws.Cells.Item(29, 11).Value = 11999
$ws.Cells.Item(29, 13).Value = -11709

$ws.Cells.Item(132, 8).Value = 1348.8536
$ws.Cells.Item(132, 9).Value = 1111.2963
$ws.Cells.Item(132, 10).Value = 1807
$ws.Cells.Item(132, 11).Value = 3333.8889
$ws.Cells.Item(132, 12).Value = 5421
$ws.Cells.Item(132, 13).Value = -803.8888999999999
$ws.Cells.Item(132, 14).Value = -10481
